# Actualización automática 2025-06-16 13:01:14
# Adds a new "GRANITO" column (inserted before the old column F) and three new
# trailing columns "NO RESURTIBLES", "PANELES PVC", "PANELES PU" to the
# "VENTAS POR GRUPO" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# ---------------------------------------------------------------------------
# 1. Insert a new column at F ("GRANITO"), shifting GRIFERIAS..SAL SOLUBLE
#    (and all their data) one column to the right, from F..N to G..O.
# ---------------------------------------------------------------------------
$ws.Columns("F:F").Insert()

# Copy number formats / styles from the neighboring (old) column into the new
# column F so header/data/footer rows keep matching look & feel.
$ws.Range("G1:G29").Copy()
$ws.Range("F1:F29").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Column width for the new GRANITO column (xlsx width 13 -> ColumnWidth 12.17)
$ws.Columns.Item(6).ColumnWidth = 13 - 0.83

# Header
$ws.Range("F1").Value = "GRANITO"

# Data rows 2-28 default to 0
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
}

# Footer row (row 29) summary text
$ws.Range("F29").Value = "0 de 27"

# ---------------------------------------------------------------------------
# 2. Append three brand-new trailing columns: P (NO RESURTIBLES),
#    Q (PANELES PVC), R (PANELES PU).
# ---------------------------------------------------------------------------
$ws.Range("O1:O29").Copy()
$ws.Range("P1:P29").PasteSpecial(-4122)
$ws.Range("Q1:Q29").PasteSpecial(-4122)
$ws.Range("R1:R29").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Columns.Item(16).ColumnWidth = 20 - 0.83
$ws.Columns.Item(17).ColumnWidth = 17 - 0.83
$ws.Columns.Item(18).ColumnWidth = 16 - 0.83

$ws.Range("P1").Value = "NO RESURTIBLES"
$ws.Range("Q1").Value = "PANELES PVC"
$ws.Range("R1").Value = "PANELES PU"

for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0
    $ws.Cells.Item($r, 18).Value = 0
}

$ws.Range("P29").Value = "0 de 27"
$ws.Range("Q29").Value = "0 de 27"
$ws.Range("R29").Value = "0 de 27"

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
